# Update statistics for June, 2024
#
# The "JuneRaw" sheet (hidden helper sheet) was blank (all zeros/empty);
# this fills it in with the real monthly circulation numbers for June.
# The visible "June" sheet pulls every value from JuneRaw via formulas
# (=JuneRaw!B2 etc.), and "Yearly total" sums January..December, so both
# recalculate automatically once JuneRaw has real data.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("MayRaw")
$dst = $wb.Worksheets.Item("JuneRaw")

# Row labels (column A) and the three column headers (B1:D1) are identical
# across every month's "Raw" sheet - copy them from MayRaw so the shared
# strings line up exactly instead of retyping the library names.
$dst.Range("A1:A55").Value() = $src.Range("A1:A55").Value()
$dst.Range("B1:D1").Value() = $src.Range("B1:D1").Value()

# Monthly figures for June 2024: Items owned by this library checked out
# at this library this month (B), items owned by other libraries checked
# out at this library this month (C), total circulation this month (D).
# Rows 11, 26, 40-44 are section/heading rows with no numbers (blank).
$data = @{
    2 = @(5395, 1367, 6762)
    3 = @(3687, 648, 4335)
    4 = @(9939, 1120, 11059)
    5 = @(290, 44, 334)
    6 = @(6449, 1008, 7457)
    7 = @(511, 259, 770)
    8 = @(826, 177, 1003)
    9 = @(678, 45, 723)
    10 = @(299, 7, 306)
    12 = @(101, 36, 137)
    13 = @(267, 141, 408)
    14 = @(494, 97, 591)
    15 = @(551, 82, 633)
    16 = @(411, 65, 476)
    17 = @(1933, 573, 2506)
    18 = @(180, 71, 251)
    19 = @(2172, 482, 2654)
    20 = @(13, 1, 14)
    21 = @(2427, 640, 3067)
    22 = @(170, 23, 193)
    23 = @(2352, 668, 3020)
    24 = @(8860, 1338, 10198)
    25 = @(952, 182, 1134)
    27 = @(579, 277, 856)
    28 = @(299, 61, 360)
    29 = @(1690, 394, 2084)
    30 = @(22, 44, 66)
    31 = @(594, 32, 626)
    32 = @(1813, 408, 2221)
    33 = @(1800, 385, 2185)
    34 = @(689, 186, 875)
    35 = @(8972, 762, 9734)
    36 = @(1239, 223, 1462)
    37 = @(3683, 494, 4177)
    38 = @(280, 35, 315)
    39 = @(171, 54, 225)
    45 = @(571, 91, 662)
    46 = @(1455, 445, 1900)
    47 = @(5056, 1034, 6090)
    48 = @(2743, 215, 2958)
    49 = @(1233, 504, 1737)
    50 = @(4671, 969, 5640)
    51 = @(624, 203, 827)
    52 = @(1480, 320, 1800)
    53 = @(355, 103, 458)
    54 = @(347, 28, 375)
    55 = @(401, 363, 764)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $dst.Cells.Item($row, 2).Value() = $vals[0]
    $dst.Cells.Item($row, 3).Value() = $vals[1]
    $dst.Cells.Item($row, 4).Value() = $vals[2]
}
